$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("25 Mar 2020", "0,6", "-12.5,6", "-25,6", "-37.5,6", "-30,6", "-20,6", "-32.5,6", "-17.5,6", "-27.5,6"),
    @("30 Mar 2020", "0,6", "-12.5,6", "-25,6", "-37.5,6", "-30,6", "-20,6", "-32.5,6", "-17.5,6", "-27.5,6"),
    @("9 Apr 2020",  "0,6", "-10,6",   "-22.5,6", "-35,6", "-27.5,6", "-17.5,6", "-30,6", "-15,6", "-25,6"),
    @("15 Apr 2020", "0,6", "-10,6",   "-22.5,6", "-35,6", "-27.5,6", "-17.5,6", "-30,6", "-15,6", "-25,6"),
    @("21 Apr 2020", "0,6", "-10,6",   "-22.5,6", "-35,6", "-27.5,6", "-17.5,6", "-30,6", "-15,6", "-25,6")
)

$startRow = 17
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowValues[$c]
    }
}

$ws.Range("E24").Select()
